$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.535365
$ws.Cells.Item(2,8).Value = 1.606095
$ws.Cells.Item(2,9).Value = 0.1618182173563651
$ws.Cells.Item(2,10).Value = 0.1618182173563651
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.01036033333333333
$ws.Cells.Item(2,14).Value = 0.031081
$ws.Cells.Item(2,15).Value = 0.0003369947480386084
$ws.Cells.Item(2,16).Value = 0.0003369947480386084
$ws.Cells.Item(2,17).Value = 0.005546559855
$ws.Cells.Item(2,18).Value = 0.049919038695
$ws.Cells.Item(2,19).Value = 0.00005453188938606502
$ws.Cells.Item(2,20).Value = 0.00005453188938606502

# Row 3
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.535365
$ws.Cells.Item(3,8).Value = 1.606095
$ws.Cells.Item(3,9).Value = 0.1618182173563651
$ws.Cells.Item(3,10).Value = 0.1618182173563651
$ws.Cells.Item(3,15).Value = 0.8439700329797517
$ws.Cells.Item(3,16).Value = 0.8439700329797518
$ws.Cells.Item(3,17).Value = 13.89081085386
$ws.Cells.Item(3,18).Value = 125.01729768474
$ws.Cells.Item(3,19).Value = 0.1365697262389761
$ws.Cells.Item(3,20).Value = 0.1365697262389761

# Row 4
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.535365
$ws.Cells.Item(4,8).Value = 1.606095
$ws.Cells.Item(4,9).Value = 0.1618182173563651
$ws.Cells.Item(4,10).Value = 0.1618182173563651
$ws.Cells.Item(4,15).Value = 0.1556929722722096
$ws.Cells.Item(4,16).Value = 0.1556929722722096
$ws.Cells.Item(4,17).Value = 2.562533673705
$ws.Cells.Item(4,18).Value = 23.062803063345
$ws.Cells.Item(4,19).Value = 0.02519395922800294
$ws.Cells.Item(4,20).Value = 0.02519395922800294

# Row 5
$ws.Cells.Item(5,9).Value = 0.6224306076670297
$ws.Cells.Item(5,10).Value = 0.6224306076670296
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.01036033333333333
$ws.Cells.Item(5,14).Value = 0.031081
$ws.Cells.Item(5,15).Value = 0.0003369947480386084
$ws.Cells.Item(5,16).Value = 0.0003369947480386084
$ws.Cells.Item(5,17).Value = 0.02133473398366667
$ws.Cells.Item(5,18).Value = 0.192012605853
$ws.Cells.Item(5,19).Value = 0.0002097558458022686
$ws.Cells.Item(5,20).Value = 0.0002097558458022686

# Row 6
$ws.Cells.Item(6,9).Value = 0.6224306076670297
$ws.Cells.Item(6,10).Value = 0.6224306076670296
$ws.Cells.Item(6,15).Value = 0.8439700329797517
$ws.Cells.Item(6,16).Value = 0.8439700329797518
$ws.Cells.Item(6,19).Value = 0.52531278048035
$ws.Cells.Item(6,20).Value = 0.5253127804803499

# Row 7
$ws.Cells.Item(7,9).Value = 0.6224306076670297
$ws.Cells.Item(7,10).Value = 0.6224306076670296
$ws.Cells.Item(7,15).Value = 0.1556929722722096
$ws.Cells.Item(7,16).Value = 0.1556929722722096
$ws.Cells.Item(7,19).Value = 0.09690807134087744
$ws.Cells.Item(7,20).Value = 0.09690807134087744

# Row 8
$ws.Cells.Item(8,9).Value = 0.2157511749766052
$ws.Cells.Item(8,10).Value = 0.2157511749766052
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.01036033333333333
$ws.Cells.Item(8,14).Value = 0.031081
$ws.Cells.Item(8,15).Value = 0.0003369947480386084
$ws.Cells.Item(8,16).Value = 0.0003369947480386084
$ws.Cells.Item(8,17).Value = 0.007395192119555556
$ws.Cells.Item(8,18).Value = 0.06655672907599999
$ws.Cells.Item(8,19).Value = 0.00007270701285027478
$ws.Cells.Item(8,20).Value = 0.00007270701285027479

# Row 9
$ws.Cells.Item(9,9).Value = 0.2157511749766052
$ws.Cells.Item(9,10).Value = 0.2157511749766052
$ws.Cells.Item(9,15).Value = 0.8439700329797517
$ws.Cells.Item(9,16).Value = 0.8439700329797518
$ws.Cells.Item(9,19).Value = 0.1820875262604257
$ws.Cells.Item(9,20).Value = 0.1820875262604257

# Row 10
$ws.Cells.Item(10,9).Value = 0.2157511749766052
$ws.Cells.Item(10,10).Value = 0.2157511749766052
$ws.Cells.Item(10,15).Value = 0.1556929722722096
$ws.Cells.Item(10,16).Value = 0.1556929722722096
$ws.Cells.Item(10,19).Value = 0.03359094170332923
$ws.Cells.Item(10,20).Value = 0.03359094170332925

